# Swap the presentation's colour theme: the deck currently uses the
# "Integral" theme colours on the (single) slide master design; the
# target state uses the classic "Office Theme" colours instead (while
# the font scheme and effect/format scheme - already identical between
# the two named themes in this deck - stay as-is).
#
# PowerPoint's ThemeColorScheme exposes exactly the 12 theme colour
# slots (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink) in this fixed
# index order, each as a settable .RGB (VBA-style 0xBBGGRR integer).

function ToRGB([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

# Target ("Office Theme") colours, in ThemeColorScheme index order.
$officeThemeColors = @(
    (ToRGB 0x00 0x00 0x00),  # 1  dk1       000000
    (ToRGB 0xFF 0xFF 0xFF),  # 2  lt1       FFFFFF
    (ToRGB 0x44 0x54 0x6A),  # 3  dk2       44546A
    (ToRGB 0xE7 0xE6 0xE6),  # 4  lt2       E7E6E6
    (ToRGB 0x5B 0x9B 0xD5),  # 5  accent1   5B9BD5
    (ToRGB 0xED 0x7D 0x31),  # 6  accent2   ED7D31
    (ToRGB 0xA5 0xA5 0xA5),  # 7  accent3   A5A5A5
    (ToRGB 0xFF 0xC0 0x00),  # 8  accent4   FFC000
    (ToRGB 0x44 0x72 0xC4),  # 9  accent5   4472C4
    (ToRGB 0x70 0xAD 0x47),  # 10 accent6   70AD47
    (ToRGB 0x05 0x63 0xC1),  # 11 hlink     0563C1
    (ToRGB 0x95 0x4F 0x72)   # 12 folHlink  954F72
)

$p  = $ppt.ActivePresentation
$cs = $p.SlideMaster.Theme.ThemeColorScheme

for ($i = 1; $i -le $cs.Count; $i++) {
    $cs.Item($i).RGB = $officeThemeColors[$i - 1]
}
